# Refresh the cryptos price/volume snapshot (commit: "Updated cryptos list on
# Tue Jun 18 14:13:19 UTC 2024 with GitHub Actions"). Rows 2-51 of Sheet1 get new
# Price (D) / Volume(1h) (E) readings pulled from the source feed; two rows also
# swapped rank (WrappedEther/WrappedBTC traded places at #15/#16), so their Coin
# name, Link and Price/Volume columns are rewritten too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Price column D holds numeric-looking strings ("7.00", "0.800", "25.80", ...)
    # that the scraper stores as plain text with digits/zeros preserved exactly.
    # A direct .Value assignment would let Excel re-interpret them as numbers and
    # silently drop significant trailing zeros, so momentarily force text entry and
    # then restore the default (unstyled) cell format, matching the source file.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "64.488.25"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.384.52"
$ws.Range("E3").Value = "  -3.50%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue "D5" "577.79"
$ws.Range("E5").Value = "  -3.70%  "
Set-TextValue "D6" "133.88"
$ws.Range("E6").Value = "  -6.45%  "
$ws.Range("D8").Value = "3.383.67"
$ws.Range("E8").Value = "  -3.49%  "
Set-TextValue "D9" "0.488"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("E10").Value = "  -10.84%  "
Set-TextValue "D11" "6.96"
$ws.Range("E11").Value = "  -11.01%  "
$ws.Range("E12").Value = "  -8.46%  "
$ws.Range("D13").Value = "3.961.48"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("E14").Value = "  -11.77%  "
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.569.86"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.369.87"
$ws.Range("E17").Value = "  -3.85%  "
Set-TextValue "D18" "25.78"
$ws.Range("E18").Value = "  -9.40%  "
Set-TextValue "D19" "9.44"
$ws.Range("E19").Value = "  -14.01%  "
$ws.Range("E20").Value = "  -7.43%  "
Set-TextValue "D21" "13.25"
$ws.Range("E21").Value = "  -7.26%  "
Set-TextValue "D22" "376.15"
$ws.Range("E22").Value = "  -9.40%  "
Set-TextValue "D23" "0.542"
$ws.Range("E23").Value = "  -9.20%  "
$ws.Range("E24").Value = "  +0.05%  "
Set-TextValue "D25" "71.48"
$ws.Range("E25").Value = "  -7.39%  "
$ws.Range("D26").Value = "3.520.10"
$ws.Range("E26").Value = "  -3.48%  "
Set-TextValue "D27" "0.0000101"
$ws.Range("E27").Value = "  -11.25%  "
$ws.Range("E28").Value = "  +0.09%  "
Set-TextValue "D29" "2.16"
$ws.Range("E29").Value = "  -11.19%  "
$ws.Range("E30").Value = "  -10.35%  "
Set-TextValue "D31" "7.85"
$ws.Range("E31").Value = "  -11.12%  "
$ws.Range("D32").Value = "3.395.47"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("E33").Value = "  -0.05%  "
Set-TextValue "D34" "22.81"
$ws.Range("E34").Value = "  -6.00%  "
Set-TextValue "D35" "0.139"
$ws.Range("E35").Value = "  -8.35%  "
Set-TextValue "D36" "168.27"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("E37").Value = "  -12.31%  "
$ws.Range("E38").Value = "  -13.80%  "
$ws.Range("E39").Value = "  -9.21%  "
Set-TextValue "D40" "4.55"
$ws.Range("E40").Value = "  -13.41%  "
Set-TextValue "D41" "0.0739"
$ws.Range("E41").Value = "  -9.28%  "
Set-TextValue "D42" "0.800"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("E43").Value = "  +0.28%  "
Set-TextValue "D44" "41.66"
$ws.Range("E44").Value = "  -7.63%  "
Set-TextValue "D45" "4.26"
$ws.Range("E45").Value = "  -15.94%  "
Set-TextValue "D46" "1.56"
$ws.Range("E46").Value = "  -11.81%  "
Set-TextValue "D47" "1.07"
$ws.Range("E47").Value = "  -0.50%  "
Set-TextValue "D48" "21.74"
$ws.Range("E48").Value = "  -6.38%  "
Set-TextValue "D49" "6.36"
$ws.Range("E49").Value = "  -9.31%  "
$ws.Range("D50").Value = "2.142.19"
$ws.Range("E50").Value = "  -8.09%  "
Set-TextValue "D51" "1.96"
$ws.Range("E51").Value = "  -16.52%  "
